# Update latest output (run 173)
# Applies refreshed optimisation-result numbers to the "Schedule" and
# "Detailed" sheets, matching the latest solver run.

$wb = $excel.ActiveWorkbook

# --- Schedule sheet: refreshed cost / unit-cost figures for rows 3-4 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("E3").Value = 722.596641
$wsSchedule.Range("F3").Value = 27.30901893424036
$wsSchedule.Range("E4").Value = 377.408538
$wsSchedule.Range("F4").Value = 11.0937253968254

# --- Detailed sheet: refreshed price forecast/history values ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B37").Value = -8.114129999999999
$wsDetailed.Range("B38").Value = -1.63893

$wsDetailed.Range("B39").Value = -1.22557
$wsDetailed.Range("C39").Value = "historical"

$wsDetailed.Range("B40").Value = 57.31
$wsDetailed.Range("C40").Value = "historical"

$wsDetailed.Range("B43").Value = 56.99272
$wsDetailed.Range("B44").Value = 50.32042
$wsDetailed.Range("B45").Value = 50.14997
$wsDetailed.Range("B46").Value = 29.68691
$wsDetailed.Range("B47").Value = 57.03877

$wsDetailed.Range("B51").Value = 57.06

$wsDetailed.Range("B54").Value = 48.14048
$wsDetailed.Range("B55").Value = 49.36549
$wsDetailed.Range("B56").Value = 56.97994

$wsDetailed.Range("B58").Value = 57.06

$wsDetailed.Range("B61").Value = 58.05476
$wsDetailed.Range("B62").Value = 59.09159

$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B65").Value = 34.27959
$wsDetailed.Range("B66").Value = 22.07
$wsDetailed.Range("B67").Value = 22.07

$wsDetailed.Range("B69").Value = 1.01974
$wsDetailed.Range("B70").Value = 0.7
$wsDetailed.Range("B71").Value = 3.73016
$wsDetailed.Range("B72").Value = 22.8184
$wsDetailed.Range("B73").Value = 23.15941
$wsDetailed.Range("B74").Value = 23.32323
$wsDetailed.Range("B75").Value = 26.88071
$wsDetailed.Range("B76").Value = 36.06
$wsDetailed.Range("B77").Value = 24.40325
$wsDetailed.Range("B78").Value = 29.97319
$wsDetailed.Range("B79").Value = 36.06
$wsDetailed.Range("B80").Value = 36.06
$wsDetailed.Range("B81").Value = 7.898
$wsDetailed.Range("B82").Value = 0.6899999999999999
$wsDetailed.Range("B83").Value = -4.59814

$wsDetailed.Range("B85").Value = -5.33831
$wsDetailed.Range("B86").Value = 12.11395
$wsDetailed.Range("B87").Value = 36.00655

$wsDetailed.Range("B90").Value = 57.06

$wsDetailed.Range("B92").Value = 56.21612
$wsDetailed.Range("B93").Value = 56.24168

Write-Output "Applied run 173 updates to Schedule and Detailed sheets."
